$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    0.04231489763667018,
    1.228061995268202,
    0.5137754236260815,
    0.4028038717171413,
    0.5336388157440486,
    -1.102517691576566,
    0.1737007515684039,
    0.8389029408811082,
    -0.6954484448595206,
    -0.2465870357053012,
    -0.1975260465718366,
    0.4425040297996861,
    -0.2720610750631522,
    -0.1065518669046048,
    -0.1895682054566924,
    1.157000698704573,
    -0.4886691766355519,
    1.10624937372658,
    -0.6446211617534254,
    -0.6387305113048862,
    0.3668428211138005,
    -0.4578680368388337,
    0.4181606776922825
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2

    # Shift existing B:J values (10 columns, cols 2..10) right into C:K (cols 3..11)
    for ($col = 11; $col -ge 3; $col--) {
        $srcVal = $ws.Cells.Item($row, $col - 1).Value()
        $ws.Cells.Item($row, $col).Value = $srcVal
    }

    # Insert the new value into column B
    $ws.Cells.Item($row, 2).Value = $newValues[$i]
}
